$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format price column as Text up front so numeric-looking strings
# (e.g. "1.0000", "239.06") are preserved exactly as typed, matching
# the source data which stores these as literal strings.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 47/48: RenderToken and Aptos swap places with updated data
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# Update Price (D) and Volume(1h) (E) columns
$ws.Range("D2").Value = "29.356.85"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.839.86"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "239.06"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "0.6260"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.07388"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "0.2887"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "24.78"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.839.00"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "4.961"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "0.6666"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "81.47"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "6.242"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "29.309.85"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "234.40"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "157.29"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").Value = "8.458"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "0.1340"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").Value = "17.29"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "0.07231"
$ws.Range("E28").Value = "  +10.62%  "
$ws.Range("D29").Value = "1.493"
$ws.Range("E29").Value = "  +4.18%  "
$ws.Range("D30").Value = "1.483"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "4.026"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "1.815"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").Value = "0.7072"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "2.586"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "0.01831"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").Value = "2.786"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "1.233.33"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").Value = "6.761"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").Value = "0.9513"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "1.985.36"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "101.03"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "65.09"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("D47").Value = "6.943"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "1.689"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").Value = "8.904"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "0.1128"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "0.3872"
$ws.Range("E51").Value = "  -2.36%  "
